$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 11.84718356105861
$ws.Range("C2").Value = 4.607969009783151
$ws.Range("D2").Value = 6.008113633237622
$ws.Range("E2").Value = 11.42787351104395
$ws.Range("G2").Value = 43.96017917940841
$ws.Range("H2").Value = 17.76681824632242
$ws.Range("K2").Value = 11.2041859675038
$ws.Range("M2").Value = 14.85762540428166
$ws.Range("N2").Value = 20.95012636998845
$ws.Range("B3").Value = 11.59417118009497
$ws.Range("C3").Value = 4.378386650634695
$ws.Range("D3").Value = 5.893884132112629
$ws.Range("E3").Value = 11.21331066153708
$ws.Range("G3").Value = 43.66106301643755
$ws.Range("H3").Value = 17.77027154000032
$ws.Range("K3").Value = 11.03704569432956
$ws.Range("M3").Value = 14.70109058664244
$ws.Range("N3").Value = 20.9952902600068
$ws.Range("B4").Value = 11.43962585447226
$ws.Range("C4").Value = 4.239167254690751
$ws.Range("D4").Value = 5.824472146972985
$ws.Range("E4").Value = 11.08326795785615
$ws.Range("G4").Value = 43.48824968428494
$ws.Range("H4").Value = 17.77557050980009
$ws.Range("K4").Value = 10.9363879636816
$ws.Range("M4").Value = 14.60836868834186
$ws.Range("N4").Value = 21.02490336762466
$ws.Range("B5").Value = 11.37695002007681
$ws.Range("C5").Value = 4.181410050280204
$ws.Range("D5").Value = 5.796411029342147
$ws.Range("E5").Value = 11.03077860664354
$ws.Range("G5").Value = 43.4206094098477
$ws.Range("H5").Value = 17.77852863282153
$ws.Range("K5").Value = 10.89591868280674
$ws.Range("M5").Value = 14.57147759095519
$ws.Range("N5").Value = 21.03744425920668
$ws.Range("B6").Value = 11.36656398124762
$ws.Range("C6").Value = 4.171760500965948
$ws.Range("D6").Value = 5.791766323156175
$ws.Range("E6").Value = 11.02209547167543
$ws.Range("G6").Value = 43.40954737108147
$ws.Range("H6").Value = 17.77906804743442
$ws.Range("K6").Value = 10.88923356845601
$ws.Range("M6").Value = 14.56540699980219
$ws.Range("N6").Value = 21.03955525148174
$ws.Range("B7").Value = 11.43877922740148
$ws.Range("C7").Value = 4.238392340892828
$ws.Range("D7").Value = 5.82409273843981
$ws.Range("E7").Value = 11.08255792587158
$ws.Range("G7").Value = 43.48732612818817
$ws.Range("H7").Value = 17.77560717098609
$ws.Range("K7").Value = 10.9358398837873
$ws.Range("M7").Value = 14.6078674901939
$ws.Range("N7").Value = 21.02507058211328
$ws.Range("B8").Value = 11.75983727072639
$ws.Range("C8").Value = 4.523191926076421
$ws.Range("D8").Value = 5.968602132476984
$ws.Range("E8").Value = 11.35358549540856
$ws.Range("G8").Value = 43.85482783637823
$ws.Range("H8").Value = 17.7673489755486
$ws.Range("K8").Value = 11.14618093646257
$ws.Range("M8").Value = 14.80297361806484
$ws.Range("N8").Value = 20.96530820150575
$ws.Range("B9").Value = 12.39152986490939
$ws.Range("C9").Value = 5.156641304608186
$ws.Range("D9").Value = 6.255894674870112
$ws.Range("E9").Value = 11.8951972288346
$ws.Range("G9").Value = 44.65887393134859
$ws.Range("H9").Value = 17.77638960343712
$ws.Range("K9").Value = 11.57182237032753
$ws.Range("M9").Value = 15.21058949250362
$ws.Range("N9").Value = 20.86304991218057
$ws.Range("B10").Value = 12.85138721278451
$ws.Range("C10").Value = 5.579553303646978
$ws.Range("D10").Value = 6.466936844336139
$ws.Range("E10").Value = 12.29486835931426
$ws.Range("G10").Value = 45.2965378365577
$ws.Range("H10").Value = 17.79841937752349
$ws.Range("K10").Value = 11.88928615961119
$ws.Range("M10").Value = 15.52267961914675
$ws.Range("N10").Value = 20.79702639448914
$ws.Range("B11").Value = 13.05857666525306
$ws.Range("C11").Value = 5.761174528986896
$ws.Range("D11").Value = 6.562451769428831
$ws.Range("E11").Value = 12.47616967752937
$ws.Range("G11").Value = 45.59590550802315
$ws.Range("H11").Value = 17.81177671163715
$ws.Range("K11").Value = 12.03405428156417
$ws.Range("M11").Value = 15.6668180392527
$ws.Range("N11").Value = 20.76896830761623
$ws.Range("B12").Value = 13.13666185008296
$ws.Range("C12").Value = 5.828405224071227
$ws.Range("D12").Value = 6.598512131344357
$ws.Range("E12").Value = 12.54467905768266
$ws.Range("G12").Value = 45.71052264976407
$ws.Range("H12").Value = 17.8173133155316
$ws.Range("K12").Value = 12.08886933456682
$ws.Range("M12").Value = 15.72166213403682
$ws.Range("N12").Value = 20.75862777349555
$ws.Range("B13").Value = 13.1198626942375
$ws.Range("C13").Value = 5.813994570129085
$ws.Range("D13").Value = 6.590751333609429
$ws.Range("E13").Value = 12.52993191818117
$ws.Range("G13").Value = 45.6857834129285
$ws.Range("H13").Value = 17.8160996519146
$ws.Range("K13").Value = 12.07706505198339
$ws.Range("M13").Value = 15.70983961268815
$ws.Range("N13").Value = 20.76084213699654
$ws.Range("B14").Value = 13.06500872361941
$ws.Range("C14").Value = 5.766736591735604
$ws.Range("D14").Value = 6.565420877785488
$ws.Range("E14").Value = 12.48180930179378
$ws.Range("G14").Value = 45.60531049935599
$ws.Range("H14").Value = 17.81222263036448
$ws.Range("K14").Value = 12.038564374526
$ws.Range("M14").Value = 15.67132512870421
$ws.Range("N14").Value = 20.76811188320893
$ws.Range("B15").Value = 13.03135807718369
$ws.Range("C15").Value = 5.737588558692434
$ws.Range("D15").Value = 6.549889922636947
$ws.Range("E15").Value = 12.45231178261439
$ws.Range("G15").Value = 45.5561792015525
$ws.Range("H15").Value = 17.80991011144129
$ws.Range("K15").Value = 12.01497924539724
$ws.Range("M15").Value = 15.64776654532603
$ws.Range("N15").Value = 20.77260186505677
$ws.Range("B16").Value = 12.83779921651114
$ws.Range("C16").Value = 5.567467527134764
$ws.Range("D16").Value = 6.460681521218969
$ws.Range("E16").Value = 12.28300337692739
$ws.Range("G16").Value = 45.27715338568597
$ws.Range("H16").Value = 17.7976134741847
$ws.Range("K16").Value = 11.87982747301146
$ws.Range("M16").Value = 15.51329919795424
$ws.Range("N16").Value = 20.79889984717345
$ws.Range("B17").Value = 12.71848217136109
$ws.Range("C17").Value = 5.460348409684499
$ws.Range("D17").Value = 6.405801601326915
$ws.Range("E17").Value = 12.17895498498312
$ws.Range("G17").Value = 45.10830160516966
$ws.Range("H17").Value = 17.79092352570903
$ws.Range("K17").Value = 11.7969659933763
$ws.Range("M17").Value = 15.43132755604359
$ws.Range("N17").Value = 20.8155391739391
$ws.Range("B18").Value = 12.64967049157888
$ws.Range("C18").Value = 5.39772299014702
$ws.Range("D18").Value = 6.374192275107752
$ws.Range("E18").Value = 12.11906510595813
$ws.Range("G18").Value = 45.0120622049988
$ws.Range("H18").Value = 17.78738972852109
$ws.Range("K18").Value = 11.74934268224518
$ws.Range("M18").Value = 15.3843865850256
$ws.Range("N18").Value = 20.82529568038481
$ws.Range("B19").Value = 12.62634317218403
$ws.Range("C19").Value = 5.376344998943984
$ws.Range("D19").Value = 6.363483580673702
$ws.Range("E19").Value = 12.09878211344157
$ws.Range("G19").Value = 44.97963076750494
$ws.Range("H19").Value = 17.78624722135257
$ws.Range("K19").Value = 11.73322624037219
$ws.Range("M19").Value = 15.36853023530957
$ws.Range("N19").Value = 20.82863100500348
$ws.Range("B20").Value = 12.73120331689662
$ws.Range("C20").Value = 5.471856284348101
$ws.Range("D20").Value = 6.411648487215143
$ws.Range("E20").Value = 12.19003619951855
$ws.Range("G20").Value = 45.12618571462488
$ws.Range("H20").Value = 17.7916031832912
$ws.Range("K20").Value = 11.80578338200069
$ws.Range("M20").Value = 15.44003251848137
$ws.Range("N20").Value = 20.81374863643737
$ws.Range("B21").Value = 13.08113143833428
$ws.Range("C21").Value = 5.780659320697326
$ws.Range("D21").Value = 6.572864305588581
$ws.Range("E21").Value = 12.49594858647837
$ws.Range("G21").Value = 45.62891400954753
$ws.Range("H21").Value = 17.8133484315038
$ws.Range("K21").Value = 12.04987354475578
$ws.Range("M21").Value = 15.68263105028903
$ws.Range("N21").Value = 20.76596886034761
$ws.Range("B22").Value = 13.30762009690683
$ws.Range("C22").Value = 5.973474424388752
$ws.Range("D22").Value = 6.67757588760881
$ws.Range("E22").Value = 12.69500101455035
$ws.Range("G22").Value = 45.96473990712865
$ws.Range("H22").Value = 17.8303482881415
$ws.Range("K22").Value = 12.2093437940431
$ws.Range("M22").Value = 15.84268752266658
$ws.Range("N22").Value = 20.73640016727259
$ws.Range("B23").Value = 13.18696774278426
$ws.Range("C23").Value = 5.87138832693219
$ws.Range("D23").Value = 6.621761307342642
$ws.Range("E23").Value = 12.58886631586474
$ws.Range("G23").Value = 45.78486668926012
$ws.Range("H23").Value = 17.82102053105088
$ws.Range("K23").Value = 12.1242548063798
$ws.Range("M23").Value = 15.75714106730004
$ws.Range("N23").Value = 20.75202972348823
$ws.Range("B24").Value = 12.72545274925275
$ws.Range("C24").Value = 5.466656817633562
$ws.Range("D24").Value = 6.409005288235901
$ws.Range("E24").Value = 12.18502659834087
$ws.Range("G24").Value = 45.11809769900693
$ws.Range("H24").Value = 17.7912949369783
$ws.Range("K24").Value = 11.80179698999812
$ws.Range("M24").Value = 15.43609642299415
$ws.Range("N24").Value = 20.81455754548334
$ws.Range("B25").Value = 12.22100100601765
$ws.Range("C25").Value = 4.991345469042464
$ws.Range("D25").Value = 6.178009603297473
$ws.Range("E25").Value = 11.74805609059154
$ws.Range("G25").Value = 44.43281778965409
$ws.Range("H25").Value = 17.77124058204236
$ws.Range("K25").Value = 11.45560596218655
$ws.Range("M25").Value = 15.09790853213653
$ws.Range("N25").Value = 20.88911455479964
